$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (character-width units); the host quantizes to the nearest
# 1/6-character pixel step, so 368.833333 is the closest input that lands the
# stored width next to the target 369.6.
$ws.Columns.Item(3).ColumnWidth = 368.833333

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = '2025-04-28 02:31:33'
$ws.Range("C47").Value = 'James Davis moved New Battery from Ford from floor space 2 to floor space 1.
Now James Davis is Frustrated.
'
$ws.Range("C47").WrapText = $true

$ws.Range("A48").Value = 47
$ws.Range("B48").Value = '2025-04-28 02:35:14'
$ws.Range("C48").Value = 'James Davis received New Battery from Ford from Suppliers Old Reliable.
New Battery from Ford''s state was New.
Thus James Davis carried out the following actions:
Update Battery Status, .
Now James Davis is Tired.
'
$ws.Range("C48").WrapText = $true

$ws.Range("A49").Value = 48
$ws.Range("B49").Value = '2025-04-28 02:48:44'
$ws.Range("C49").Value = 'James Davis moved battery 7 from floor space 1 to floor space 3.
Now James Davis is Happy.
'
$ws.Range("C49").WrapText = $true

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = '2025-04-28 03:03:06'
$ws.Range("C50").Value = 'James Davis took picture of New Battery from Ford.
Now James Davis is Tired, feeling that the task was Tiring.
'
$ws.Range("C50").WrapText = $true

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = '2025-04-28 03:05:54'
$ws.Range("C51").Value = 'James Davis moved New Battery from Ford from floor space 2 to floor space 2.
Now James Davis is Confident, feeling that the task was Challenging.
'
$ws.Range("C51").WrapText = $true

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = '2025-04-28 03:08:10'
$ws.Range("C52").Value = 'James Davis received Its brand new from ford from Suppliers Battery New.
Its brand new from ford''s state was Old.
Thus James Davis carried out the following actions:
Update Battery Status, Diagnostic Analysis, Disassembly, Repair, Re-assembly, .
Now James Davis is Tired, feeling that the task was Tiring.
'
$ws.Range("C52").WrapText = $true

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = '2025-04-28 03:10:04'
$ws.Range("C53").Value = 'James Davis added Telsla Battery 4 to the database.
Serial Number is 573432019330921.
Part Number is 3322.
Item Type is 3.
Location is floor space 1.
Now James Davis is Frustrated, feeling that the task was Challenging.
'
$ws.Range("C53").WrapText = $true
